$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match-by-match rows scraped for Glenn Maxwell (Kings XI Punjab), appended
# below the existing data (rows 2-12), extending the table through row 23.
$newRows = @(
    @(" Abu Dhabi", " October 30 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Glenn Maxwell ", "6", "6", "1", "0", "100.00"),
    @(" Dubai (DSC)", " October 04 2020", "Super Kings won by 10 wickets (with 14 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "Glenn Maxwell ", "11", "7", "1", "0", "157.14"),
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Kings XI Punjab", "Mumbai Indians", "Glenn Maxwell ", "11", "18", "0", "0", "61.11"),
    @(" Dubai (DSC)", " September 24 2020", "Kings XI won by 97 runs", "Kings XI Punjab", "Royal Challengers Bangalore", "Glenn Maxwell ", "5", "6", "0", "0", "83.33"),
    @(" Dubai (DSC)", " October 20 2020", "Kings XI won by 5 wickets (with 6 balls remaining)", "Kings XI Punjab", "Delhi Capitals", "Glenn Maxwell ", "32", "24", "3", "0", "133.33"),
    @(" Abu Dhabi", " October 10 2020", "KKR won by 2 runs", "Kings XI Punjab", "Kolkata Knight Riders", "Glenn Maxwell ", "10", "5", "2", "0", "200.00"),
    @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Kings XI Punjab", "Delhi Capitals", "Glenn Maxwell ", "1", "4", "0", "0", "25.00"),
    @(" Sharjah", " September 27 2020", "Royals won by 4 wickets (with 3 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Glenn Maxwell ", "13", "9", "2", "0", "144.44"),
    @(" Dubai (DSC)", " October 18 2020", "Match tied (Kings XI won the one-over eliminator)", "Kings XI Punjab", "Mumbai Indians", "Glenn Maxwell ", "0", "2", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Glenn Maxwell ", "12", "13", "0", "0", "92.30"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Glenn Maxwell ", "7", "12", "0", "0", "58.33")
)

$startRow = 13
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $col = $c + 1
        # Leading apostrophe forces text (prevents Excel from coercing
        # numeric-looking strings such as '6' or '100.00' into numbers),
        # matching the source data which stores every column as text.
        $ws.Cells.Item($r, $col).Value = "'" + $rowValues[$c]
    }
}
